$d = $word.ActiveDocument

# Edit 1: merge "Men det finns många fler " + "jämförelse-operatorer." into one run's text
$d.Content.Find.Execute(
    "Men det finns många fler jämförelse-operatorer.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Men det finns många fler jämförelse-operatorer.",
    2
)

# Edit 2: extend "Dessutom kommer inte programmet köras" with additional sentence
$d.Content.Find.Execute(
    "Dessutom kommer inte programmet köras",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Dessutom kommer inte programmet köras. Dessutom kan man inte tilldela en variabel dess värde  i ett if statement. ",
    2
)
